# Course Xpath's Updated and Loop Count Increased
#
# The automation loop that stamps the "PortfolioCourse" / "AssignmentName"
# header pair on the STAGE sheet (row 2, columns M:P) ran a few more
# iterations, so the generated course/assignment identifiers rolled over to
# new values. Update the four header cells to the newly generated names —
# this naturally grows the shared-string table while the cells keep their
# existing (border + centered) look.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STAGE")

$ws.Range("M2").Value = "PortfolioCourse43630"
$ws.Range("N2").Value = "AssignmentName43630"
$ws.Range("O2").Value = "PortfolioCourse62144"
$ws.Range("P2").Value = "AssignmentName62144"
